$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 264; everything from row 264 downward
# (264-311) shifts down to 266-313.
$ws.Rows.Item(264).Insert()
$ws.Rows.Item(264).Insert()

# New row 264
$ws.Range("A264").Value = 7
$ws.Range("B264").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C264").Value = "Ñuble"
$ws.Range("D264").Value = 44644
$ws.Range("E264").Value = 16
$ws.Range("F264").Value = 100114001
$ws.Range("G264").Value = "Papa"
$ws.Range("H264").Value = "Asterix"
$ws.Range("I264").Value = "1a (cosecha lavada)"
$ws.Range("J264").Value = 100
$ws.Range("K264").Value = 8500
$ws.Range("L264").Value = 9000
$ws.Range("M264").Value = 8750
$ws.Range("N264").Value = "$/malla 25 kilos"
$ws.Range("O264").Value = "Región Metropolitana"
$ws.Range("P264").Value = 350
$ws.Range("Q264").Value = 25
$ws.Range("R264").Value = "Hortaliza"

# New row 265
$ws.Range("A265").Value = 7
$ws.Range("B265").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C265").Value = "Ñuble"
$ws.Range("D265").Value = 44644
$ws.Range("E265").Value = 16
$ws.Range("F265").Value = 100114001
$ws.Range("G265").Value = "Papa"
$ws.Range("H265").Value = "Patagonia"
$ws.Range("I265").Value = "1a (cosecha)"
$ws.Range("J265").Value = 120
$ws.Range("K265").Value = 7000
$ws.Range("L265").Value = 7500
$ws.Range("M265").Value = 7250
$ws.Range("N265").Value = "$/saco 25 kilos"
$ws.Range("O265").Value = "Provincia de Diguillín"
$ws.Range("P265").Value = 290
$ws.Range("Q265").Value = 25
$ws.Range("R265").Value = "Hortaliza"
